$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 50

# Column C holds a numeric-looking value ("25") that must stay stored as
# text (matching the rest of the "Value" column), so force the cell to
# Text format before writing it -- otherwise Excel auto-converts it to a
# number.
$ws.Range("C$row").NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($row, 2).Value = "temperature"
$ws.Cells.Item($row, 3).Value = "25"
$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "N/A"
$ws.Cells.Item($row, 6).Value = "N/A"
